$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 390, pushing existing rows 390-448 down to 392-450
$ws.Rows.Item(390).Resize(2).Insert()

# Row 390
$ws.Range("A390").Value = 7
$ws.Range("B390").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C390").Value = 'Ñuble'
$ws.Range("D390").Value = 45034
$ws.Range("E390").Value = 16
$ws.Range("F390").Value = 100112002
$ws.Range("G390").Value = 'Pimiento'
$ws.Range("H390").Value = 'Zafiro rojo'
$ws.Range("I390").Value = 'Primera'
$ws.Range("J390").Value = 60
$ws.Range("K390").Value = 17000
$ws.Range("L390").Value = 17000
$ws.Range("M390").Value = 17000
$ws.Range("N390").Value = '$/caja 15 kilos'
$ws.Range("O390").Value = 'Región de Arica y Parinacota'
$ws.Range("P390").Value = 1133
$ws.Range("Q390").Value = 15
$ws.Range("R390").Value = 'Hortaliza'

# Row 391
$ws.Range("A391").Value = 7
$ws.Range("B391").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C391").Value = 'Ñuble'
$ws.Range("D391").Value = 45034
$ws.Range("E391").Value = 16
$ws.Range("F391").Value = 100112002
$ws.Range("G391").Value = 'Pimiento'
$ws.Range("H391").Value = 'Zafiro verde'
$ws.Range("I391").Value = 'Primera'
$ws.Range("J391").Value = 60
$ws.Range("K391").Value = 12000
$ws.Range("L391").Value = 12000
$ws.Range("M391").Value = 12000
$ws.Range("N391").Value = '$/caja 15 kilos'
$ws.Range("O391").Value = 'Región de Arica y Parinacota'
$ws.Range("P391").Value = 800
$ws.Range("Q391").Value = 15
$ws.Range("R391").Value = 'Hortaliza'
